$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)
$seq = $s.TimeLine.MainSequence
for ($i = $seq.Count; $i -ge 1; $i--) {
    $seq.Item($i).Delete()
}
